# Stage the new record's values on a temporary helper sheet so we can use
# Copy/PasteSpecial to get them onto the protected "Minuta de registro"
# sheet (direct Range.Value / .Formula writes are rejected there because
# the sheet is protected and we don't have the unlock password).
$wb = $excel.ActiveWorkbook
$stage = $wb.Worksheets.Add()
$ws = $wb.Worksheets.Item("Minuta de registro")

$stage.Range("A1").Value = "6WF2MN3"
$stage.Range("B1").NumberFormat = "@"
$stage.Range("B1").Value = "46243"
$stage.Range("C1").Value = "Andrés Felipe Pérez"
$stage.Range("D1").Value = "Dell Latitude 5420"
$stage.Range("E1").Value = "PC Laptop"
$stage.Range("G1").Value = "Ingreso"
$stage.Range("H1").Value = "Andrés Felipe Pérez"
$stage.Range("I1").Value = "ANDRESFELIPE.PEREZ"

$stage.Range("A1:E1").Copy()
$ws.Range("A2").PasteSpecial()

$stage.Range("G1:I1").Copy()
$ws.Range("G2").PasteSpecial()

# F2 (registration date/time) is a genuine number, which the protected
# sheet refuses via paste, so unlock just that cell long enough to write
# the serial date value, then relock it.
$ws.Range("F2").Locked = $false
$ws.Range("F2").Value = 45027.627129629633
$ws.Range("F2").Locked = $true

$stage.Delete()

$ws.Range("A1:J1").Select()
